# phieumuonrs.xlsx - "PHIẾU MƯỢN TÀI LIỆU" (borrow slip) gets regenerated
# with a different borrower/book, and the title cell (B12) gets wrap-text
# turned on so long titles no longer get clipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mã số thẻ: (library card number)
$ws.Range("B8").Value = "2033207524"

# Phí cọc: (deposit fee)
$ws.Range("B9").Value = "250000 đồng"

# Tên tài liệu / Tác giả / KH xếp giá (title / author / shelf code)
$ws.Range("B12").Value = "Lotharingia: A Personal History Of Europe's Lost Country"
$ws.Range("C12").Value = "Simon Winder"
$ws.Range("D12").Value = "B-300"

# Turn on wrap text for the (now longer) title cell so it fits the row.
$ws.Range("B12").WrapText = $true

# Ngày ... tháng ... năm ... (date line)
$ws.Range("C16").Value = "TP. Hồ Chí Minh, Ngày 17 tháng 1 năm 2021."

# Matches the last-edited cell recorded in the saved selection.
$ws.Range("B12").Select()
